$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '332.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.05%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '45.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.54%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.563'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.25%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08337'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.19%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.047'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '6.17%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9809'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.14%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.82%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1942'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '5.44%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '10.30'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-13.65%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1010'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.25%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04632'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.42%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.71%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001269'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.58%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006028'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.97%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.365'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.32%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.436'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2.83%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.631'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.42%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3348'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-3.65%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1385'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.25%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2490'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.03%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04114'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.18%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001302'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4.44%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004423'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.13%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001278'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '7.09%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003739'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.53%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02820'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '10.97%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05768'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '5.89%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007646'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.67%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1428'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.10%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007555'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.26%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.49%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008032'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.45%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007200'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.87%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.33%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005799'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.21%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003488'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-27.54%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.33%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.33%'
